$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text formatting (values like "6.640" or
# "0.9980" would otherwise be auto-converted to numbers and lose trailing
# zeros), matching the source workbook where these columns are stored as text.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.382.55"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.864.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9982"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.43"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.18%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9982"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07766"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3087"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.177"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "93.47"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.849.03"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6952"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.640"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008374"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.359.79"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.14"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.096.60"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.82"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9982"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.589"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9987"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1521"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.936"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.78"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.46"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.539"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.260"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.219"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.205"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05161"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7939"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.18%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.155"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.693"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.337.48"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +8.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.732"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9571"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +6.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.071"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +9.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "107.85"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9980"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.814"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.996.96"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.24"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5191"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.781"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.77%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.025"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.27%  "
